# Apply the "season record" columns (Wins, Losses, Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row: copy the style from the existing header cell AC1 so the new
# headers match the look of the rest of row 1 (bold, centered, bordered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 85  # column AD
    $ws.Cells.Item($r, 31).Value = 77  # column AE
    $ws.Cells.Item($r, 32).Value = 0   # column AF
}
